$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3274.375
$ws.Range("J19").Value = 3219.8
$ws.Range("L19").Value = 3219.8
$ws.Range("N19").Value = -3569.8
$ws.Range("H53").Value = 2642.8462
$ws.Range("I53").Value = 2580.1538
$ws.Range("J53").Value = 2705.5386
$ws.Range("K53").Value = 2580.1538
$ws.Range("L53").Value = 2705.5386
$ws.Range("M53").Value = -1943.1538
$ws.Range("N53").Value = -3979.5386
$ws.Range("H55").Value = 229
$ws.Range("J55").Value = 73.75
$ws.Range("L55").Value = 73.75
$ws.Range("N55").Value = -501.75
$ws.Range("H131").Value = 4001.8333
$ws.Range("I131").Value = 3441.25
$ws.Range("K131").Value = 10323.75
$ws.Range("M131").Value = -5283.75
$ws.Range("H132").Value = 2248.9443
$ws.Range("I132").Value = 2279.375
$ws.Range("J132").Value = 2005.5
$ws.Range("K132").Value = 6838.125
$ws.Range("L132").Value = 6016.5
$ws.Range("M132").Value = -4308.125
$ws.Range("N132").Value = -11076.5
$ws.Range("H135").Value = 500715.25
$ws.Range("I135").Value = 625745
$ws.Range("K135").Value = 5631705
$ws.Range("M135").Value = -5629170
$ws.Range("H138").Value = 4109.74
$ws.Range("I138").Value = 1885.1666
$ws.Range("J138").Value = 5361.0625
$ws.Range("K138").Value = 5655.4998
$ws.Range("L138").Value = 16083.1875
$ws.Range("M138").Value = -515.4997999999996
$ws.Range("N138").Value = -26363.1875
$ws.Range("H141").Value = 13335490
$ws.Range("I141").Value = 13891079
$ws.Range("K141").Value = 41673237
$ws.Range("M141").Value = -41668057

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 320.2857
$ws.Range("I5").Value = 323.33334
$ws.Range("K5").Value = 323.33334
$ws.Range("M5").Value = -211.33334
$ws.Range("H61").Value = 3470.5334
$ws.Range("I61").Value = 1560.0156
$ws.Range("K61").Value = 1560.0156
$ws.Range("M61").Value = -1348.0156
$ws.Range("H74").Value = 38249.566
$ws.Range("I74").Value = 51308.516
$ws.Range("J74").Value = 5099.923
$ws.Range("K74").Value = 51308.516
$ws.Range("L74").Value = 5099.923
$ws.Range("M74").Value = -50434.516
$ws.Range("N74").Value = -6847.923
$ws.Range("H77").Value = 38249.566
$ws.Range("I77").Value = 51308.516
$ws.Range("J77").Value = 5099.923
$ws.Range("K77").Value = 256542.58
$ws.Range("L77").Value = 25499.615
$ws.Range("M77").Value = -252174.58
$ws.Range("N77").Value = -34235.615
$ws.Range("H102").Value = 962.2857
$ws.Range("I102").Value = 956
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 956
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = 666
$ws.Range("N102").Value = -4244
$ws.Range("H136").Value = 3470.5334
$ws.Range("I136").Value = 1560.0156
$ws.Range("K136").Value = 4680.0468
$ws.Range("M136").Value = -2130.0468

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 320.2857
$ws.Range("I4").Value = 323.33334
$ws.Range("K4").Value = 323.33334
$ws.Range("M4").Value = -208.33334
$ws.Range("H86").Value = 38504676
$ws.Range("I86").Value = 66496.69
$ws.Range("J86").Value = 100005760
$ws.Range("K86").Value = 66496.69
$ws.Range("L86").Value = 100005760
$ws.Range("M86").Value = -65373.69
$ws.Range("N86").Value = -100008006
$ws.Range("H89").Value = 38504676
$ws.Range("I89").Value = 66496.69
$ws.Range("J89").Value = 100005760
$ws.Range("K89").Value = 332483.45
$ws.Range("L89").Value = 500028800
$ws.Range("M89").Value = -326867.45
$ws.Range("N89").Value = -500040032
$ws.Range("H99").Value = 3248927
$ws.Range("I99").Value = 2037.9565
$ws.Range("K99").Value = 2037.9565
$ws.Range("M99").Value = -539.9565
$ws.Range("H134").Value = 4152.3115
$ws.Range("I134").Value = 1089.8889
$ws.Range("K134").Value = 3269.6667
$ws.Range("M134").Value = -734.6666999999998

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 64
$ws.Range("J7").Value = 100
$ws.Range("L7").Value = 100
$ws.Range("N7").Value = -326
$ws.Range("H31").Value = 7220
$ws.Range("I31").Value = 2552.238
$ws.Range("J31").Value = 11481.869
$ws.Range("K31").Value = 2552.238
$ws.Range("L31").Value = 11481.869
$ws.Range("M31").Value = -2257.238
$ws.Range("N31").Value = -12071.869
$ws.Range("H34").Value = 7220
$ws.Range("I34").Value = 2552.238
$ws.Range("J34").Value = 11481.869
$ws.Range("K34").Value = 2552.238
$ws.Range("L34").Value = 11481.869
$ws.Range("M34").Value = -2350.238
$ws.Range("N34").Value = -11885.869
$ws.Range("H58").Value = 9808908
$ws.Range("I58").Value = 16668060
$ws.Range("J58").Value = 10118.619
$ws.Range("K58").Value = 16668060
$ws.Range("L58").Value = 10118.619
$ws.Range("M58").Value = -16667857
$ws.Range("N58").Value = -10524.619
$ws.Range("H88").Value = 27583.334
$ws.Range("J88").Value = 27583.334
$ws.Range("L88").Value = 27583.334
$ws.Range("N88").Value = -28395.334
$ws.Range("H91").Value = 27583.334
$ws.Range("J91").Value = 27583.334
$ws.Range("L91").Value = 27583.334
$ws.Range("N91").Value = -30391.334
$ws.Range("H134").Value = 7789.2964
$ws.Range("I134").Value = 1827.2858
$ws.Range("K134").Value = 5481.857400000001
$ws.Range("M134").Value = -2946.857400000001
$ws.Range("H136").Value = 9808908
$ws.Range("I136").Value = 16668060
$ws.Range("J136").Value = 10118.619
$ws.Range("K136").Value = 50004180
$ws.Range("L136").Value = 30355.857
$ws.Range("M136").Value = -50001630
$ws.Range("N136").Value = -35455.857

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H114").Value = 701.7222
$ws.Range("I114").Value = 540.5714
$ws.Range("J114").Value = 804.2727
$ws.Range("K114").Value = 1621.7142
$ws.Range("L114").Value = 2412.8181
$ws.Range("M114").Value = 1632.2858
$ws.Range("N114").Value = -8920.8181
$ws.Range("H116").Value = 2438.8
$ws.Range("I116").Value = 2438.8
$ws.Range("K116").Value = 7316.400000000001
$ws.Range("M116").Value = -3874.400000000001

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 727682.0600000001
$ws.Range("J107").Value = 833.3333
$ws.Range("L107").Value = 833.3333
$ws.Range("N107").Value = -4673.3333
$ws.Range("H122").Value = 4832228
$ws.Range("I122").Value = 5574724.5
$ws.Range("K122").Value = 16724173.5
$ws.Range("M122").Value = -16721723.5
$ws.Range("H132").Value = 4951.579
$ws.Range("I132").Value = 1970
$ws.Range("J132").Value = 13300
$ws.Range("K132").Value = 5910
$ws.Range("L132").Value = 39900
$ws.Range("M132").Value = -3380
$ws.Range("N132").Value = -44960

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1892.7333
$ws.Range("I22").Value = 1025
$ws.Range("K22").Value = 1025
$ws.Range("M22").Value = -730
$ws.Range("H27").Value = 1892.7333
$ws.Range("I27").Value = 1025
$ws.Range("K27").Value = 1025
$ws.Range("M27").Value = -918
$ws.Range("H46").Value = 1726305.6
$ws.Range("I46").Value = 5747641.5
$ws.Range("J46").Value = 2876
$ws.Range("K46").Value = 5747641.5
$ws.Range("L46").Value = 2876
$ws.Range("M46").Value = -5747453.5
$ws.Range("N46").Value = -3252
$ws.Range("H122").Value = 3616.5574
$ws.Range("I122").Value = 2792.26
$ws.Range("J122").Value = 7363.364
$ws.Range("K122").Value = 8376.780000000001
$ws.Range("L122").Value = 22090.092
$ws.Range("M122").Value = -5926.780000000001
$ws.Range("N122").Value = -26990.092
$ws.Range("H132").Value = 8202414
$ws.Range("I132").Value = 15627847
$ws.Range("K132").Value = 46883541
$ws.Range("M132").Value = -46881011

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 3990200.5
$ws.Range("J5").Value = 6333667.5
$ws.Range("L5").Value = 6333667.5
$ws.Range("N5").Value = -6333891.5
$ws.Range("H88").Value = 59900
$ws.Range("I88").Value = 59900
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 59900
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -59494
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 59900
$ws.Range("I91").Value = 59900
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 59900
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -58496
$ws.Range("N91").ClearContents()
$ws.Range("H122").Value = 111310.19
$ws.Range("I122").Value = 135158.06
$ws.Range("J122").Value = 9105
$ws.Range("K122").Value = 405474.18
$ws.Range("L122").Value = 27315
$ws.Range("M122").Value = -403024.18
$ws.Range("N122").Value = -32215
$ws.Range("H126").Value = 1400.6
$ws.Range("I126").Value = 832.5454999999999
$ws.Range("K126").Value = 2497.6365
$ws.Range("M126").Value = -27.63649999999961
$ws.Range("H132").Value = 9267409
$ws.Range("J132").Value = 19174.375
$ws.Range("L132").Value = 57523.125
$ws.Range("N132").Value = -62583.125
$ws.Range("H136").Value = 16969826
$ws.Range("I136").Value = 25642152
$ws.Range("K136").Value = 76926456
$ws.Range("M136").Value = -76923906
